$d = $word.ActiveDocument

# "Becoming an UI Testing Rock Star" -> "Becoming a UI Testing Rock Star"
$d.Content.Find.Execute("Becoming an UI Testing Rock Star", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Becoming a UI Testing Rock Star", 2)
